$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 50015.953
$ws.Range("I11").Value = 50015.953
$ws.Range("K11").Value = 50015.953
$ws.Range("M11").Value = -49875.953
$ws.Range("H43").Value = 2620
$ws.Range("I43").Value = 1200
$ws.Range("K43").Value = 1200
$ws.Range("M43").Value = -1131
$ws.Range("H75").Value = 62693.5
$ws.Range("J75").Value = 64332.332
$ws.Range("L75").Value = 64332.332
$ws.Range("N75").Value = -66204.33199999999
$ws.Range("H78").Value = 62693.5
$ws.Range("J78").Value = 64332.332
$ws.Range("L78").Value = 192996.996
$ws.Range("N78").Value = -202356.996
$ws.Range("H100").Value = 4759.154
$ws.Range("I100").Value = 5089.4546
$ws.Range("K100").Value = 5089.4546
$ws.Range("M100").Value = -4548.4546
$ws.Range("H101").Value = 515.13336
$ws.Range("J101").Value = 912.8
$ws.Range("L101").Value = 2738.4
$ws.Range("N101").Value = -5982.4
$ws.Range("H113").Value = 38465040
$ws.Range("I113").Value = 62503256
$ws.Range("K113").Value = 62503256
$ws.Range("M113").Value = -62500002
$ws.Range("H132").Value = 4506.549
$ws.Range("I132").Value = 4056.8838
$ws.Range("K132").Value = 12170.6514
$ws.Range("M132").Value = -9640.651400000001
$ws.Range("H137").Value = 2138.7
$ws.Range("I137").Value = 1780.8334
$ws.Range("J137").Value = 2377.2778
$ws.Range("K137").Value = 5342.5002
$ws.Range("L137").Value = 7131.8334
$ws.Range("M137").Value = -2792.5002
$ws.Range("N137").Value = -12231.8334
$ws.Range("H138").Value = 2136.7026
$ws.Range("I138").Value = 2530.7144
$ws.Range("J138").Value = 2044.7667
$ws.Range("K138").Value = 7592.1432
$ws.Range("L138").Value = 6134.300099999999
$ws.Range("M138").Value = -2452.1432
$ws.Range("N138").Value = -16414.3001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 81215.39999999999
$ws.Range("I102").Value = 54764.117
$ws.Range("J102").Value = 253148.75
$ws.Range("K102").Value = 54764.117
$ws.Range("L102").Value = 253148.75
$ws.Range("M102").Value = -53142.117
$ws.Range("N102").Value = -256392.75
$ws.Range("H122").Value = 2280.0645
$ws.Range("J122").Value = 4099.7
$ws.Range("L122").Value = 12299.1
$ws.Range("N122").Value = -17199.1
$ws.Range("H133").Value = 49354.855
$ws.Range("J133").Value = 49617.832
$ws.Range("L133").Value = 49617.832
$ws.Range("N133").Value = -54677.832

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 10871267
$ws.Range("I94").Value = 13158886
$ws.Range("J94").Value = 5075.5
$ws.Range("K94").Value = 13158886
$ws.Range("L94").Value = 5075.5
$ws.Range("M94").Value = -13158435
$ws.Range("N94").Value = -5977.5
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 869.7143
$ws.Range("J22").Value = 817.8
$ws.Range("L22").Value = 817.8
$ws.Range("N22").Value = -1517.8
$ws.Range("H31").Value = 3455.9673
$ws.Range("I31").Value = 2817.8
$ws.Range("J31").Value = 3512.9465
$ws.Range("K31").Value = 2817.8
$ws.Range("L31").Value = 3512.9465
$ws.Range("M31").Value = -2522.8
$ws.Range("N31").Value = -4102.9465
$ws.Range("H34").Value = 3455.9673
$ws.Range("I34").Value = 2817.8
$ws.Range("J34").Value = 3512.9465
$ws.Range("K34").Value = 2817.8
$ws.Range("L34").Value = 3512.9465
$ws.Range("M34").Value = -2615.8
$ws.Range("N34").Value = -3916.9465
$ws.Range("H88").Value = 13832.444
$ws.Range("J88").Value = 13832.444
$ws.Range("L88").Value = 13832.444
$ws.Range("N88").Value = -14644.444
$ws.Range("H91").Value = 13832.444
$ws.Range("J91").Value = 13832.444
$ws.Range("L91").Value = 13832.444
$ws.Range("N91").Value = -16640.444
$ws.Range("H105").Value = 421.2857
$ws.Range("I105").Value = 423.07693
$ws.Range("K105").Value = 423.07693
$ws.Range("M105").Value = 1323.92307

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 4333.3335
$ws.Range("I8").Value = 4333.3335
$ws.Range("K8").Value = 13000.0005
$ws.Range("M8").Value = -12861.0005
$ws.Range("H63").Value = 19114.6
$ws.Range("J63").Value = 8250
$ws.Range("L63").Value = 24750
$ws.Range("N63").Value = -26248
$ws.Range("H66").Value = 19114.6
$ws.Range("J66").Value = 8250
$ws.Range("L66").Value = 74250
$ws.Range("N66").Value = -81738
$ws.Range("H99").Value = 4248.5835
$ws.Range("I99").Value = 1245.75
$ws.Range("K99").Value = 3737.25
$ws.Range("M99").Value = -1491.25
$ws.Range("H120").Value = 20000
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 0
$ws.Range("M120").ClearContents()

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 473.7
$ws.Range("I97").Value = 454.9565
$ws.Range("J97").Value = 535.2857
$ws.Range("K97").Value = 454.9565
$ws.Range("L97").Value = 535.2857
$ws.Range("M97").Value = 41.04349999999999
$ws.Range("N97").Value = -1527.2857
$ws.Range("H122").Value = 2169.1304
$ws.Range("I122").Value = 1905.6154
$ws.Range("J122").Value = 2511.7
$ws.Range("K122").Value = 5716.8462
$ws.Range("L122").Value = 7535.099999999999
$ws.Range("M122").Value = -3266.8462
$ws.Range("N122").Value = -12435.1
$ws.Range("H132").Value = 3958.7036
$ws.Range("J132").Value = 2561.3845
$ws.Range("L132").Value = 7684.1535
$ws.Range("N132").Value = -12744.1535
$ws.Range("H139").Value = 147998.6
$ws.Range("J139").Value = 193331.67
$ws.Range("L139").Value = 193331.67
$ws.Range("N139").Value = -203611.67

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H34").Value = 5000
$ws.Range("I34").Value = 5000
$ws.Range("K34").Value = 5000
$ws.Range("M34").Value = -4828
$ws.Range("H46").Value = 2830.423
$ws.Range("J46").Value = 3135.238
$ws.Range("L46").Value = 3135.238
$ws.Range("N46").Value = -3511.238
$ws.Range("H61").Value = 1881.0769
$ws.Range("I61").Value = 2361.889
$ws.Range("K61").Value = 2361.889
$ws.Range("M61").Value = -2159.889
$ws.Range("H93").Value = 111114450
$ws.Range("I93").Value = 111114450
$ws.Range("K93").Value = 111114450
$ws.Range("M93").Value = -111113202
$ws.Range("H100").Value = 3051.8823
$ws.Range("I100").Value = 3058.8
$ws.Range("K100").Value = 3058.8
$ws.Range("M100").Value = -2517.8
$ws.Range("H113").Value = 1881.0769
$ws.Range("I113").Value = 2361.889
$ws.Range("K113").Value = 2361.889
$ws.Range("M113").Value = -191.8890000000001
$ws.Range("H136").Value = 6302
$ws.Range("J136").Value = 3448.2307
$ws.Range("L136").Value = 10344.6921
$ws.Range("N136").Value = -15444.6921

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2171820.2
$ws.Range("I62").Value = 7940841
$ws.Range("K62").Value = 7940841
$ws.Range("M62").Value = -7940217
$ws.Range("H65").Value = 2171820.2
$ws.Range("I65").Value = 7940841
$ws.Range("K65").Value = 39704205
$ws.Range("M65").Value = -39701085
$ws.Range("H81").Value = 13338068
$ws.Range("I81").Value = 3265.125
$ws.Range("J81").Value = 28577842
$ws.Range("K81").Value = 6530.25
$ws.Range("L81").Value = 57155684
$ws.Range("M81").Value = -5469.25
$ws.Range("N81").Value = -57157806
$ws.Range("H84").Value = 13338068
$ws.Range("I84").Value = 3265.125
$ws.Range("J84").Value = 28577842
$ws.Range("K84").Value = 32651.25
$ws.Range("L84").Value = 285778420
$ws.Range("M84").Value = -27347.25
$ws.Range("N84").Value = -285789028
$ws.Range("H132").Value = 6508.1304
$ws.Range("I132").Value = 7029.375
$ws.Range("K132").Value = 21088.125
$ws.Range("M132").Value = -18558.125
$ws.Range("H136").Value = 3575.0952
$ws.Range("I136").Value = 3201
$ws.Range("K136").Value = 9603
$ws.Range("M136").Value = -7053
